# Reviewed Process_mapping.R and Analyse_MC_Tool_data.py
# Fill in the "basic purpose / non-obvious details / suggestions" columns
# for the "Process mining" (Process_mapping.R) and "Main Folder"
# (Analyse_MC_Tool_Data) rows on Sheet1, and update the sheet's active
# selection to reflect where the review ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21: Process mining / Process_mapping.R -----------------------
$ws.Range("D21").Value = "Script for performing process mining on model output data. This is useful for validating modelling pathway and for understanding bottlenecks in model"
$ws.Range("E21").Value = "Reads in distribution_outputs.csv  data from deterministic simulations, only.`nIs run after running model `nReads event log data"
$ws.Range("F21").Value = "There's a whole host of graphs produced in this script, both for process maps and for general EA. `nThink it would be worth using this as the basis for all outputs and going away from other png files, below."
$ws.Range("D21:F21").WrapText = $true
$ws.Rows.Item(21).RowHeight = 72

# --- Row 22: Main Folder / Analyse_MC_Tool_Data ------------------------
$ws.Range("D22").Value = "Analyse Monte-Carlo data"
$ws.Range("E22").Value = "Reads in monte_carlo_table.csv (simple data on discharge time and admission time) and computes  statistics across runs looking at variability of 4 hr percentages."
$ws.Range("F22").Value = "Might want to expand script to not only look at 4 hr percentages, but a range of other metrics. Suggest adjusting monte_carlo_table to take full event logs, as in distribution_outputs.csv not only discharge and admission times."
$ws.Range("E22:F22").WrapText = $true
$ws.Rows.Item(22).RowHeight = 57.6

# --- Reflect where the reviewer left the selection ---------------------
$ws.Range("F23").Select()
try {
    $win = $excel.ActiveWindow
    $win.ScrollRow = 13
    $win.ScrollColumn = 2
} catch {
    # Best-effort only; scroll position isn't critical data.
}
